$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# above the existing row for this market/product (currently row 94), so
# insert a new row there - this shifts all the old rows 94-123 down to 95-124.
$ws.Rows(94).Insert()

# Populate the newly inserted row 94 with the new weekly record.
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44782
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = "Fruta"
$ws.Cells.Item(94, 7).Value = 100107
$ws.Cells.Item(94, 8).Value = "Otros"
$ws.Cells.Item(94, 9).Value = 100107002
$ws.Cells.Item(94, 10).Value = "Chirimoya"
$ws.Cells.Item(94, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 35
$ws.Cells.Item(94, 14).Value = 4000
$ws.Cells.Item(94, 15).Value = 4000
$ws.Cells.Item(94, 16).Value = 4000
$ws.Cells.Item(94, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(94, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 19).Value = 4000
$ws.Cells.Item(94, 20).Value = 1

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(94, 4).NumberFormat = $ws.Cells.Item(95, 4).NumberFormat
